$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds text-formatted numbers (e.g. "41.760.07", "303.16").
# Force text format first so Excel doesn't auto-coerce them to floats/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.760.07"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "2.267.53"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "303.16"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "92.18"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "32.39"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "53.34"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "2.618.05"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "14.23"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "2.269.50"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").Value = "41.654.92"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  +5.82%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "67.01"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "239.43"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").Value = "23.92"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -4.47%  "
$ws.Range("D31").Value = "35.26"
$ws.Range("E31").Value = "  +7.00%  "
$ws.Range("D32").Value = "160.15"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").Value = "3.01"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "16.92"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "2.006.69"
$ws.Range("D44").Value = "19.53"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").Value = "10.34"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "52.37"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +1.33%  "
